$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row-by-row updates of Price (D) and Volume(1h) (E) columns, plus the
# MXToken/ARBITRUM row swap (B/C/D/E on rows 34-35), per the scraper run.
# Numeric-looking price strings are entered with a leading apostrophe so
# Excel stores them as text (matching the source "Price" column, which is
# always text, e.g. to preserve trailing zeros / multi-dot big numbers)
# instead of auto-converting them to real numbers.

$ws.Range("D2").Value = "26.117.32"
$ws.Range("E2").Value = "  -0.68%  "
$ws.Range("D3").Value = "1.655.74"
$ws.Range("E3").Value = "  -0.78%  "
$ws.Range("D5").Value = "'218.16"
$ws.Range("E5").Value = "  -0.37%  "
$ws.Range("D6").Value = "'0.5285"
$ws.Range("E6").Value = "  +0.82%  "
$ws.Range("D8").Value = "'0.2605"
$ws.Range("E8").Value = "  -2.66%  "
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("E10").Value = "  -3.14%  "
$ws.Range("E11").Value = "  +0.27%  "
$ws.Range("D12").Value = "'4.495"
$ws.Range("E12").Value = "  +0.81%  "
$ws.Range("D13").Value = "1.657.72"
$ws.Range("E13").Value = "  -0.89%  "
$ws.Range("D14").Value = "'0.5471"
$ws.Range("E14").Value = "  -1.22%  "
$ws.Range("D15").Value = "0.0₅8162"
$ws.Range("E15").Value = "  -1.57%  "
$ws.Range("D16").Value = "'65.48"
$ws.Range("E16").Value = "  +0.55%  "
$ws.Range("D17").Value = "26.131.39"
$ws.Range("E17").Value = "  -0.72%  "
$ws.Range("D18").Value = "'1.003"
$ws.Range("E18").Value = "  -0.35%  "
$ws.Range("D19").Value = "'4.564"
$ws.Range("E19").Value = "  -2.69%  "
$ws.Range("D20").Value = "'192.87"
$ws.Range("E20").Value = "  -1.23%  "
$ws.Range("E21").Value = "  -0.81%  "
$ws.Range("D22").Value = "'6.022"
$ws.Range("E22").Value = "  -1.14%  "
$ws.Range("D23").Value = "'1.004"
$ws.Range("E23").Value = "  -0.39%  "
$ws.Range("D24").Value = "'141.84"
$ws.Range("E24").Value = "  +1.46%  "
$ws.Range("D25").Value = "'0.1250"
$ws.Range("E25").Value = "  +0.73%  "
$ws.Range("D26").Value = "'7.267"
$ws.Range("E26").Value = "  +0.60%  "
$ws.Range("D27").Value = "'16.21"
$ws.Range("E27").Value = "  -0.28%  "
$ws.Range("D28").Value = "'1.439"
$ws.Range("E28").Value = "  +1.61%  "
$ws.Range("D29").Value = "'0.05939"
$ws.Range("E29").Value = "  -4.33%  "
$ws.Range("D30").Value = "'1.278"
$ws.Range("E30").Value = "  -0.45%  "
$ws.Range("E31").Value = "  -2.23%  "
$ws.Range("D32").Value = "'3.246"
$ws.Range("E32").Value = "  -1.76%  "
$ws.Range("E33").Value = "  -3.79%  "
$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").Value = "'0.9492"
$ws.Range("E34").Value = "  -2.58%  "
$ws.Range("B35").Value = "MXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D35").Value = "'2.791"
$ws.Range("E35").Value = "  -0.14%  "
$ws.Range("D36").Value = "'2.411"
$ws.Range("E36").Value = "  -0.66%  "
$ws.Range("D37").Value = "'0.5665"
$ws.Range("E37").Value = "  -1.90%  "
$ws.Range("D38").Value = "'0.01610"
$ws.Range("E38").Value = "  -0.02%  "
$ws.Range("D39").Value = "'5.816"
$ws.Range("E39").Value = "  -3.42%  "
$ws.Range("D40").Value = "'0.8474"
$ws.Range("E41").Value = "  -0.21%  "
$ws.Range("D42").Value = "'102.57"
$ws.Range("E42").Value = "  +2.17%  "
$ws.Range("D43").Value = "1.022.15"
$ws.Range("E43").Value = "  -0.30%  "
$ws.Range("D44").Value = "1.800.13"
$ws.Range("E44").Value = "  -0.60%  "
$ws.Range("D45").Value = "'57.16"
$ws.Range("E45").Value = "  -1.13%  "
$ws.Range("E46").Value = "  +0.09%  "
$ws.Range("D47").Value = "'0.4289"
$ws.Range("E47").Value = "  +1.65%  "
$ws.Range("E48").Value = "  -1.15%  "
$ws.Range("D49").Value = "'0.05151"
$ws.Range("D50").Value = "'7.815"
$ws.Range("E50").Value = "  -3.63%  "
$ws.Range("D51").Value = "'0.09694"
$ws.Range("E51").Value = "  -1.04%  "
